$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.6545652718822623;     C = 1.626987699542094;    D = 0.1496068669990043;  E = 0.5333859586016987;  G = 2.964545797025059 }
    3 = @{ B = 0.00009552326474482342; C = 0.002658071450198252; D = 3.223369029078222;   E = 0.5333859586016987;  G = 3.759508582394863 }
    4 = @{ B = 0.003078177322033415;   C = 0.3048912486333797;   D = 0.1496068669990043;  E = 0.5333859586016987;  G = 0.9909622515561161 }
    5 = @{ B = 0.003078177322033415;   C = 1.626987699542094;    D = 18.71679738969934;   E = 13.86384647080068;   G = 34.21070973736415 }
    6 = @{ B = 1.445647641019636;      C = 1.626987699542094;    D = 3.223369029078222;   E = 13.86384647080068;   G = 20.15985084044064 }
    7 = @{ B = 3.272327238179451;      C = 0.3048912486333797;   D = 189.6080260415259;   E = 0.5333859586016987;  G = 193.7186304869404 }
    8 = @{ B = 1.445647641019636;      C = 1.626987699542094;    D = 0.7210945179870265;  E = 13.86384647080068;   G = 17.65757632934944 }
    9 = @{ B = 1.445647641019636;      C = 1.626987699542094;    D = 0.7210945179870265;  E = 0.5333859586016987;  G = 4.327115817150455 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
